# edit.ps1 - apply "Updated partial final report" changes
# Rebuilds the "What we did well" / "What we could've done better" / "What
# have we learned" sections with the new subsections described by the diff:
#   - "Sound" heading/bookmark is replaced by "Camera & User controls",
#     followed by new "Bullet Library" subsection, followed by a
#     (re-created) "Sound" subsection carrying the original bookmark name
#     and body paragraph.
#   - A new blank paragraph is added after the "Sound" body paragraph.
#   - A new "Scoring" subsection is added after "User Interfaces & CEGUI".
#   - A new paragraph is added after the "What have we learned" heading.

$d = $word.ActiveDocument

function New-DocxFragment([string]$bodyXml) {
    return '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Replace a single paragraph's whole range (incl. its end-of-paragraph mark)
# with one or more paragraphs of raw OOXML. Used for in-place rewrites.
function Set-ParagraphXml($paragraph, [string]$bodyXml) {
    $paragraph.Range.InsertXML((New-DocxFragment $bodyXml))
}

# Insert one or more paragraphs of raw OOXML right after $paragraph, leaving
# $paragraph itself untouched. Implemented as: add one blank paragraph right
# after the anchor (InsertParagraphAfter — known-good boundary behaviour),
# then overwrite that new blank paragraph's full range with the real XML
# (which itself may expand into several paragraphs).
function Add-ParagraphsAfter($paragraph, [string]$bodyXml) {
    $anchorEnd = $paragraph.Range.End
    $paragraph.Range.InsertParagraphAfter()
    $newPara = $d.Range($anchorEnd, $anchorEnd).Paragraphs(1)
    Set-ParagraphXml $newPara $bodyXml
}

# --- Step 1: replace the "Sound" Heading2 paragraph (under "What did well")
# with: Camera & User controls (heading+body), Bullet Library (heading+body),
# and a freshly recreated "Sound" heading (keeping its original bookmark name).
$soundHeadingPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Heading 2" -and $p.Range.Text.TrimEnd([char]13) -eq "Sound") {
        $soundHeadingPara = $p
        break
    }
}
if ($soundHeadingPara -eq $null) { throw "Could not find the 'Sound' Heading2 paragraph" }
$block1Xml = '<w:p><w:pPr><w:pStyle w:val="Heading2"/><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:bookmarkStart w:colFirst="0" w:colLast="0" w:name="_fsgeh2y4dxj5" w:id="100"/><w:bookmarkEnd w:id="100"/><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">' + 'Camera &amp; User controls' + '</w:t></w:r></w:p>' + '<w:p><w:pPr><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">' + 'Most of the controls were based off of the previous assignment’s code which handled movement through the keys which we transitioned to using the mouse. A teammate figured out that we can use Clamp in order to restrict where both the camera and paddle move.' + '</w:t></w:r></w:p>' + '<w:p><w:pPr><w:pStyle w:val="Heading2"/><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:bookmarkStart w:colFirst="0" w:colLast="0" w:name="_jx07v9th36td" w:id="101"/><w:bookmarkEnd w:id="101"/><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">' + 'Bullet Library' + '</w:t></w:r></w:p>' + '<w:p><w:pPr><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">' + 'Integrating Bullet wasn’t too difficult and learning how to make the ball bounce within the second week set the tone for the rest of the physics based interactions.' + '</w:t></w:r></w:p>' + '<w:p><w:pPr><w:pStyle w:val="Heading2"/><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:bookmarkStart w:colFirst="0" w:colLast="0" w:name="_ykm3epup2fkv" w:id="102"/><w:bookmarkEnd w:id="102"/><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">' + 'Sound' + '</w:t></w:r></w:p>'
Set-ParagraphXml $soundHeadingPara $block1Xml

# --- Step 2: insert a new blank paragraph right after the "Sound" body
# paragraph ("Sound was easy to implement...").
$soundBodyPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Sound was easy to implement")) {
        $soundBodyPara = $p
        break
    }
}
if ($soundBodyPara -eq $null) { throw "Could not find the Sound body paragraph" }
$block2Xml = '<w:p><w:pPr><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>'
Add-ParagraphsAfter $soundBodyPara $block2Xml

# --- Step 3: insert "Scoring" heading + body paragraph right after the
# "User Interfaces & CEGUI" body paragraph ("We used CEGUI...").
$ceguiBodyPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("We used CEGUI")) {
        $ceguiBodyPara = $p
        break
    }
}
if ($ceguiBodyPara -eq $null) { throw "Could not find the CEGUI body paragraph" }
$block3Xml = '<w:p><w:pPr><w:pStyle w:val="Heading2"/><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:bookmarkStart w:colFirst="0" w:colLast="0" w:name="_asj9vv6qi5u6" w:id="103"/><w:bookmarkEnd w:id="103"/><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">' + 'Scoring' + '</w:t></w:r></w:p>' + '<w:p><w:pPr><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">' + 'It was easy to implement a count and updating the GUI to reflect that. However, we ran into a bug that caused the ball to count multiple times as it hit the wall since it would be there for longer than one frame.' + '</w:t></w:r></w:p>'
Add-ParagraphsAfter $ceguiBodyPara $block3Xml

# --- Step 4: insert a new paragraph right after the "What have we learned"
# Heading1 paragraph.
$learnedPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Heading 1" -and $p.Range.Text.TrimEnd([char]13) -eq "What have we learned") {
        $learnedPara = $p
        break
    }
}
if ($learnedPara -eq $null) { throw "Could not find the 'What have we learned' paragraph" }
$block4Xml = '<w:p><w:pPr><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">' + 'Using hard-coded values makes adjusting objects tedious.' + '</w:t></w:r></w:p>'
Add-ParagraphsAfter $learnedPara $block4Xml

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
